$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.322.10"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "'2.642.97"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'601.84"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "'146.30"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "'0.365"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "'27.29"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'3.115.83"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'63.209.78"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "'2.628.32"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'11.43"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").Value = "'341.61"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").Value = "'6.90"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  -3.50%  "
$ws.Range("D24").Value = "'66.43"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").Value = "'1.65"
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("D26").Value = "'8.77"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("D27").Value = "'552.25"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Value = "'2.03"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("D34").Value = "'0.0₃0805"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'5.24"
$ws.Range("E35").Value = "  +6.55%  "
$ws.Range("D36").Value = "'166.16"
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.405"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("D39").Value = "'19.03"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "'1.88"
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "'168.31"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.74"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'22.46"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").Value = "'0.0575"
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").Value = "'0.624"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'0.0244"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'18.75"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "'1.77"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "'11.22"
$ws.Range("E51").Value = "  -1.09%  "
